# Auto-generated edit script: updates cryptocurrency price/volume data
# and rotates the coin list (rows 7-17) per the commit's symbol-list refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain-text columns (Coin name, Link URL) - safe to assign directly,
# Excel will not try to reinterpret these as numbers.
$textUpdates = @{
    'B7' = 'KuCoinToken'
    'C7' = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
    'B8' = 'MXToken'
    'C8' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'B9' = 'LiechtensteinCryptoassetsExchange'
    'C9' = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
    'B10' = 'WazirX'
    'C10' = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
    'B11' = 'MandalaExchangeToken'
    'C11' = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
    'B12' = 'BitrueCoin'
    'C12' = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
    'B13' = 'BitMartToken'
    'C13' = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
    'B14' = 'BitForexToken'
    'C14' = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
    'B15' = 'TigerCash'
    'C15' = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
    'B16' = 'LEO'
    'C16' = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
    'B17' = 'GateToken'
    'C17' = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
}

foreach ($key in $textUpdates.Keys) {
    $ws.Range($key).Value = $textUpdates[$key]
}

# Numeric-looking text columns (Price, Volume%) - these must stay as literal
# text (e.g. '308.48', '1.64%') rather than being auto-converted by Excel into
# numbers/percentages. Force the cell to Text format before assignment, then
# clear the temporary formatting so no stray style index is left behind.
$numericTextUpdates = @{
    'D2' = '308.48'
    'E2' = '1.64%'
    'D3' = '39.02'
    'E3' = '9.30%'
    'D4' = '5.084'
    'E4' = '0.93%'
    'D5' = '0.08185'
    'E5' = '3.54%'
    'D6' = '2.031'
    'E6' = '9.75%'
    'D7' = '7.910'
    'E7' = '1.61%'
    'D8' = '0.9298'
    'E8' = '1.11%'
    'D9' = '0.1420'
    'E9' = '4.50%'
    'D10' = '0.1946'
    'E10' = '2.83%'
    'D11' = '0.09305'
    'E11' = '2.75%'
    'D12' = '0.03468'
    'E12' = '-0.03%'
    'D13' = '0.09844'
    'E13' = '0.39%'
    'D14' = '0.001413'
    'E14' = '1.18%'
    'D15' = '0.005851'
    'E15' = '-4.72%'
    'D16' = '3.923'
    'E16' = '5.39%'
    'D17' = '4.180'
    'E17' = '1.99%'
    'D18' = '3.437'
    'E18' = '5.77%'
    'E19' = '0.43%'
    'D20' = '0.1303'
    'E20' = '-2.93%'
    'D21' = '4.833'
    'E21' = '-6.38%'
    'D22' = '0.2353'
    'E22' = '7.34%'
    'D23' = '0.04461'
    'E23' = '1.40%'
    'E24' = '0.55%'
    'E25' = '-9.51%'
    'D27' = '0.0001300'
    'E27' = '-0.01%'
    'D39' = '0.02115'
    'E39' = '8.74%'
    'D40' = '0.05177'
    'E40' = '1.60%'
    'D41' = '0.007486'
    'E41' = '-1.55%'
    'D42' = '0.01014'
    'E42' = '0.39%'
    'E43' = '2.11%'
    'E44' = '-1.85%'
    'D45' = '0.009679'
    'E45' = '-4.85%'
    'D46' = '0.00006298'
    'E46' = '2.41%'
    'E47' = '0.04%'
    'D49' = '0.001601'
    'E49' = '-3.50%'
    'D50' = '0.00002101'
    'E50' = '0.04%'
    'D51' = '0.0002001'
    'E51' = '0.04%'
}

foreach ($key in $numericTextUpdates.Keys) {
    $cell = $ws.Range($key)
    $cell.NumberFormat = "@"
    $cell.Value = $numericTextUpdates[$key]
    $cell.ClearFormats()
}
